# "Added body type classify"
#
# The "Body Types" sheet listed three blocks (Ectomorph Type, Mesomorph,
# Endomorph) separated by blank rows. This edit:
#   1. Renames the "Ectomorph Type" label to "Ectomorph".
#   2. Removes the blank separator rows so the three body-type blocks sit
#      directly one after another (classification list, tidied up).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Body Types")

# Rename "Ectomorph Type" -> "Ectomorph"
$ws.Range("A5").Value2 = "Ectomorph"

# Remove the blank row gaps between the three blocks (delete bottom-up so
# earlier row numbers stay valid while deleting).
$ws.Rows("20:22").Delete()
$ws.Rows("11:13").Delete()
$ws.Rows("2:4").Delete()

$ws.Range("F7").Select()
